$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Station Status" heading -> "Overall " + "Station Status"
#    Insert a new run "Overall " right before the existing "Station
#    Status" run, forcing it to stay a separate run (same final rPr:
#    sz 56 / szCs 56, automatic colour) by diverging its formatting
#    while the neighbouring text is inserted, then restoring it.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Station Status", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$headStart = $rng.Start
$ins = $d.Range($headStart, $headStart)
$ins.InsertBefore("Overall ")
$overallRng = $d.Range($headStart, $headStart + 8)
$overallRng.Font.Bold = $true
$overallRng.Font.Bold = $false

# ---------------------------------------------------------------------
# 2) "LWA1 Status" -> "LWA" + " Station" + " Status" (three separate
#    runs, all keeping the original red colour / size of the title).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("LWA1 Status", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "LWA", 2)

$part2 = $d.Range($rng2.End, $rng2.End)
$part2.InsertAfter(" Station")
$part2.Font.Color = 255

$part3 = $d.Range($part2.End, $part2.End)
$part3.InsertAfter(" Status")
$part3.Font.Color = 65280

# Restore the shared red colour last, once both new runs already have
# their own (temporarily different) formatting, so Word keeps them as
# distinct runs instead of re-coalescing them with their neighbours.
$part2.Font.Color = 1772703
$part3.Font.Color = 1772703

# ---------------------------------------------------------------------
# 3) Remove the stray "_GoBack" bookmark left over from the last edit
#    position, without disturbing the surrounding run structure.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

Write-Output "done"
